# Fruta / hortaliza, semanal
# Update weekly price-reporting rows (Fecha, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg) on the
# "Hortaliza, Agricola del Norte S.A. de Arica - Repollo" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44658
$ws.Range("I2").Value = "Segunda"
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 600
$ws.Range("L2").Value = 650
$ws.Range("M2").Value = 625
$ws.Range("P2").Value = 625

$ws.Range("D3").Value = 44253
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 800
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = 850
$ws.Range("P3").Value = 850

$ws.Range("D4").Value = 44253
$ws.Range("I4").Value = "Tercera"
$ws.Range("J4").Value = 800
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 700
$ws.Range("M4").Value = 650
$ws.Range("P4").Value = 650

$ws.Range("D5").Value = 44229
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 760
$ws.Range("K5").Value = 550
$ws.Range("L5").Value = 600
$ws.Range("M5").Value = 575
$ws.Range("P5").Value = 575

$ws.Range("D6").Value = 44874
$ws.Range("I6").Value = "Tercera"
$ws.Range("J6").Value = 1200
$ws.Range("K6").Value = 450
$ws.Range("L6").Value = 500
$ws.Range("M6").Value = 475
$ws.Range("P6").Value = 475

$ws.Range("D8").Value = 44935
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 400
$ws.Range("L8").Value = 500
$ws.Range("M8").Value = 460
$ws.Range("P8").Value = 460

$ws.Range("D9").Value = 44201
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 800
$ws.Range("L9").Value = 900
$ws.Range("M9").Value = 850
$ws.Range("P9").Value = 850

$ws.Range("D10").Value = 44544
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 600
$ws.Range("L10").Value = 650
$ws.Range("M10").Value = 625
$ws.Range("P10").Value = 625

$ws.Range("D11").Value = 44210
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 900
$ws.Range("K11").Value = 600
$ws.Range("L11").Value = 700
$ws.Range("M11").Value = 650
$ws.Range("P11").Value = 650

$ws.Range("D13").Value = 44278
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 700
$ws.Range("K13").Value = 600
$ws.Range("L13").Value = 700
$ws.Range("M13").Value = 650
$ws.Range("P13").Value = 650

$ws.Range("D14").Value = 44278
$ws.Range("I14").Value = "Tercera"
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 600
$ws.Range("M14").Value = 550
$ws.Range("P14").Value = 550

$ws.Range("D15").Value = 44174
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 800
$ws.Range("K15").Value = 450
$ws.Range("L15").Value = 500
$ws.Range("M15").Value = 475
$ws.Range("P15").Value = 475

$ws.Range("D16").Value = 44174
$ws.Range("I16").Value = "Tercera"
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 250
$ws.Range("L16").Value = 350
$ws.Range("M16").Value = 300
$ws.Range("P16").Value = 300

$ws.Range("D18").Value = 44799
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 800
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 1200
$ws.Range("M18").Value = 1100
$ws.Range("P18").Value = 1100

$ws.Range("D19").Value = 44474
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = 600
$ws.Range("L19").Value = 700
$ws.Range("M19").Value = 650
$ws.Range("P19").Value = 650

$ws.Range("D20").Value = 44245
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 800
$ws.Range("K20").Value = 850
$ws.Range("L20").Value = 900
$ws.Range("M20").Value = 875
$ws.Range("P20").Value = 875

$ws.Range("D21").Value = 44245
$ws.Range("I21").Value = "Segunda"
$ws.Range("J21").Value = 1000
$ws.Range("K21").Value = 750
$ws.Range("L21").Value = 800
$ws.Range("M21").Value = 775
$ws.Range("P21").Value = 775
